# Adds a new bulleted requirement right after the "Номерні знаки..."
# bullet in the requirements list:
#
#   Пробіг зафіксований на одному із оглядів авто повинен бути більший
#   або рівний пробігу з попереднього огляду.
#
# with the first sentence (everything but the trailing period) wrapped in
# a bookmark named "_Hlk154790260", matching the target revision.

$d = $word.ActiveDocument

$anchorText = "Номерні знаки машини повинні відповідати встановленому формату."

# Find the paragraph that carries the anchor sentence so the insertion
# point does not depend on a hard-coded paragraph index.
$anchorPara = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*$anchorText*") {
        $anchorPara = $para
        break
    }
}

if ($anchorPara -eq $null) {
    Write-Output "ERROR: anchor paragraph not found"
} else {
    # Split off a brand-new paragraph right after the anchor one. Because
    # it is split from the anchor paragraph's own mark, it automatically
    # inherits the same pPr (pStyle "a3" / numPr ilvl 0, numId 5 / rPr
    # lang uk-UA) as the surrounding bulleted list items.
    $anchorPara.Range.InsertParagraphAfter()

    $newPara = $anchorPara.Next()
    $newRange = $newPara.Range

    $sentence = "Пробіг зафіксований на одному із оглядів авто повинен бути більший або рівний пробігу з попереднього огляду"
    $fullText = $sentence + "."

    # Insert the whole sentence (plus trailing period) as one run first so
    # it picks up the paragraph's run formatting (rPr lang uk-UA).
    $newRange.InsertAfter($fullText)

    # Bookmark only the sentence itself, excluding the trailing period -
    # Word splits the run at that boundary automatically, leaving two
    # runs that both keep the inherited rPr, just like the target markup.
    $bookmarkRange = $d.Range($newRange.Start, $newRange.Start + $sentence.Length)
    $d.Bookmarks.Add("_Hlk154790260", $bookmarkRange)

    Write-Output "Inserted paragraph: $($newPara.Range.Text)"
}
